# Update default soil characteristics + fix bug for goes rainfall
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header labels (columns X, Y, Z) in row 1
# ---------------------------------------------------------------------------
$ws.Range("X1").Value = "infiltration_cm_day"
$ws.Range("Y1").Value = "Ks_cmday"
$ws.Range("Z1").Value = "Ks_mmsec"
# Match the vertical-center alignment style used by the other header cells
$ws.Range("X1:Z1").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 2. Column K (soilDepthCM) - replace the old "=J*2.54" formula with plain
#    hard-coded values
# ---------------------------------------------------------------------------
$ws.Range("K2").Value = 76
$ws.Range("K3").Value = 200
$ws.Range("K4").Value = 41
$ws.Range("K5").Value = 200
$ws.Range("K6").Value = 200
$ws.Range("K7").Value = 200
$ws.Range("K8").Value = 200
$ws.Range("K9").Value = 200
$ws.Range("K10").Value = 20
$ws.Range("K11").Value = 25

# ---------------------------------------------------------------------------
# 3. New columns X (infiltration_cm_day) and Z (Ks_mmsec) - plain values
# ---------------------------------------------------------------------------
$ws.Range("X2").Value = 7.85
$ws.Range("X3").Value = 6.25
$ws.Range("X4").Value = 6.25
$ws.Range("X5").Value = 6.25
$ws.Range("X6").Value = 8
$ws.Range("X7").Value = 6
$ws.Range("X8").Value = 3
$ws.Range("X9").Value = 6.25
$ws.Range("X10").Value = 4
$ws.Range("X11").Value = 7

$ws.Range("Z2").Value = 28
$ws.Range("Z3").Value = 48.48
$ws.Range("Z4").Value = 40.8
$ws.Range("Z5").Value = 86.88
$ws.Range("Z6").Value = 0
$ws.Range("Z7").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("Z9").Value = 48.48
$ws.Range("Z10").Value = 8.76
$ws.Range("Z11").Value = 40.8

# ---------------------------------------------------------------------------
# 4. New column Y (Ks_cmday) - plain integer values, uses the same
#    integer number format as column B (style with numFmtId 1)
# ---------------------------------------------------------------------------
$ws.Range("Y2:Y11").NumberFormat = "0"
$ws.Range("Y2").Value = 900
$ws.Range("Y3").Value = 600
$ws.Range("Y4").Value = 600
$ws.Range("Y5").Value = 300
$ws.Range("Y6").Value = 1000
$ws.Range("Y7").Value = 300
$ws.Range("Y8").Value = 200
$ws.Range("Y9").Value = 600
$ws.Range("Y10").Value = 1000
$ws.Range("Y11").Value = 300

# ---------------------------------------------------------------------------
# 5. Column W (infiltration_cmhr) - now computed from the new Ks_cmday
#    column via =ROUND(Y/24, 3). Row 2 keeps its own (non-shared) formula,
#    rows 3:11 form a shared-formula group.
# ---------------------------------------------------------------------------
$ws.Range("W2").Formula = "=ROUND(Y2/24, 3)"
$ws.Range("W3:W11").Formula = "=ROUND(Y3/24, 3)"

# ---------------------------------------------------------------------------
# 6. Best-effort column widths for the new columns (V was resized, W:Z are
#    brand new columns)
# ---------------------------------------------------------------------------
$ws.Columns.Item(22).ColumnWidth = 12
$ws.Columns.Item(23).ColumnWidth = 18.833333333333336
$ws.Columns.Item(24).ColumnWidth = 15.999999999999998
$ws.Columns.Item(25).ColumnWidth = 13.333333333333332
$ws.Columns.Item(26).ColumnWidth = 16.833333333333336

# ---------------------------------------------------------------------------
# 7. Sheet view: unfreeze/refreeze at column B instead of O, and move the
#    active selection to Z9
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("Z9").Select()
